# PM19 Tidsregistrering af Marc.xlsx - apply commit changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Rows 23-27: the workday date shifts from 2020-03-09 to 2020-03-10 ---
$ws.Range("C23").Value = 43900
$ws.Range("C24").Value = 43900
$ws.Range("C25").Value = 43900
$ws.Range("C26").Value = 43900
$ws.Range("C27").Value = 43900

# --- Row 28: new task "Lav SD0201" by Designer on 2020-03-11, 09:50-09:55 ---
$ws.Range("A28").Value = "Lav SD0201"
$ws.Range("B28").Value = "Designer"
$ws.Range("C28").Value = 43901
$ws.Range("D28").Value = 0.40972222222222227
$ws.Range("E28").Value = 0.41319444444444442

# --- Row 29: new task "lav SD0203" by Designer on 2020-03-11, 10:00-11:20 ---
$ws.Range("A29").Value = "lav SD0203"
$ws.Range("B29").Value = "Designer"
$ws.Range("C29").Value = 43901
$ws.Range("D29").Value = 0.41666666666666669
$ws.Range("E29").Value = 0.47222222222222227

# --- Row 30: new task "oprette entities til UC02" on 2020-03-11, 12:10-14:30 ---
$ws.Range("A30").Value = "oprette entities til UC02"
$ws.Range("C30").Value = 43901
$ws.Range("D30").Value = 0.50694444444444442
$ws.Range("E30").Value = 0.60416666666666663

# --- Update the saved selection to match the author's cursor position ---
[void]$ws.Range("E31").Select()
